# ---------------------------------------------------------------------------
# Applies the PlayerPerformance_5663.xlsx edit:
#   1. Adds a new worksheet "ODI Batting Extra" (sheetId 4) at the end of the
#      workbook and populates it with per-match batting-extras data.
#   2. Removes the stray empty B5 / B10 cells on the "ODI Batting" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Clean up the empty placeholder cells on "ODI Batting" (sheet2)
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B5").ClearContents()
$odiBatting.Range("B10").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" sheet as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Re-use the same header formatting (bold / bordered / centered) already used
# by the other sheets' header rows, so the new header row matches style s="1".
$odiBatting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le 6; $col++) {
    $extra.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Row data: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# $null entries mean "leave blank" (matches the source cells that carry no value).
$rows = @(
    @("4434", 7,     "4", "0", "11.87%", "NO"),
    @("4564", 7,     "1", "0", "8.44%",  "NO"),
    @("4565", $null, $null, $null, $null, "NO"),
    @("4567", 7,     $null, $null, $null, "NO"),
    @("4586", $null, $null, $null, $null, "NO"),
    @("4590", $null, $null, $null, $null, "NO"),
    @("4592", 6,     "3", "0", "12.64%", "NO"),
    @("4634", $null, $null, $null, $null, "NO"),
    @("4638", 7,     $null, $null, $null, "NO"),
    @("4641", 5,     "0", "0", "0.97%",  "NO")
)

$r = 2
foreach ($row in $rows) {
    $matchCode = $row[0]
    $battingPosition = $row[1]
    $num4 = $row[2]
    $num6 = $row[3]
    $pctRuns = $row[4]
    $manOfMatch = $row[5]

    # MATCH_CODE is stored as text even though it looks numeric.
    $aCell = $extra.Cells.Item($r, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $matchCode

    # BATTING_POSITION is a genuine number when present, otherwise left blank.
    $bCell = $extra.Cells.Item($r, 2)
    if ($null -ne $battingPosition) {
        $bCell.Value = $battingPosition
    } else {
        $bCell.Style = "Normal"
    }

    # NUM_4 / NUM_6 are text numbers when present, otherwise left blank.
    $cCell = $extra.Cells.Item($r, 3)
    if ($null -ne $num4) {
        $cCell.NumberFormat = "@"
        $cCell.Value = $num4
    } else {
        $cCell.Style = "Normal"
    }

    $dCell = $extra.Cells.Item($r, 4)
    if ($null -ne $num6) {
        $dCell.NumberFormat = "@"
        $dCell.Value = $num6
    } else {
        $dCell.Style = "Normal"
    }

    # PERCENT_RUNS_OF_TOTAL is a text percentage string when present, otherwise blank.
    $eCell = $extra.Cells.Item($r, 5)
    if ($null -ne $pctRuns) {
        $eCell.NumberFormat = "@"
        $eCell.Value = $pctRuns
    } else {
        $eCell.Style = "Normal"
    }

    # MAN_OF_MATCH is always populated text ("NO" for every row here).
    $extra.Cells.Item($r, 6).Value = $manOfMatch

    $r++
}

$extra.Range("A1").Select() | Out-Null

# Restore the originally-active sheet ("Player Info") as the active tab, since
# adding the new sheet would otherwise leave it selected/active.
$wb.Worksheets.Item(1).Activate() | Out-Null
